$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Hulk Hogan ---
$ws.Range("A4").Value = "Hulk"
$ws.Range("B4").Value = "Hogan"
$ws.Range("C4").Value = 34001
$ws.Range("D4").Value = "M"
$ws.Range("E4").Value = "Hulk@gmail.com"
$ws.Hyperlinks.Add($ws.Range("E4"), "mailto:Hulk@gmail.com")
$ws.Range("F4").Value = 7777777777

# --- Row 5: rey meystario ---
$ws.Range("A5").Value = "rey"
$ws.Range("B5").Value = "meystario"
$ws.Range("C5").Value = 35827
$ws.Range("D5").Value = "F"
$ws.Range("E5").Value = "rey@gmail.com"
$ws.Hyperlinks.Add($ws.Range("E5"), "mailto:rey@gmail.com")
$ws.Range("F5").Value = 6666666666

# --- Row 6: Ram Gope ---
$ws.Range("A6").Value = "Ram"
$ws.Range("B6").Value = "Gope"
$ws.Range("C6").Value = 32984
$ws.Range("D6").Value = "M"
$ws.Range("E6").Value = "Ram@gmail.com"
$ws.Hyperlinks.Add($ws.Range("E6"), "mailto:Ram@gmail.com")
$ws.Range("F6").Value = 9999999999

# --- Row 7: Junu majhi ---
$ws.Range("A7").Value = "Junu"
$ws.Range("B7").Value = "majhi"
$ws.Range("C7").Value = 33635
$ws.Range("D7").Value = "F"
$ws.Range("E7").Value = "junu@gmail.com"
$ws.Hyperlinks.Add($ws.Range("E7"), "mailto:junu@gmail.com")
$ws.Range("F7").Value = 5555555555

# Match the date/hyperlink display formats used by the existing rows (C2:C3, E2:E3).
# Applying these AFTER Hyperlinks.Add keeps the cells on the same pre-existing
# "Hyperlink" style (xf) rather than a second one minted by Hyperlinks.Add.
$ws.Range("C2").Copy()
$ws.Range("C4:C7").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E4:E7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B15").Select()
